$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting (bold, border, centered) from H1, then set values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-25
$data = @{
    2  = @(1, 2)
    3  = @(1, 1)
    4  = @(1, 6)
    5  = @(2, 5)
    6  = @(3, 6)
    7  = @(1, 5)
    8  = @(4, 9)
    9  = @(1, 6)
    10 = @(1, 4)
    11 = @(1, 5)
    12 = @(1, 5)
    13 = @(1, 4)
    14 = @(1, 5)
    15 = @(1, 6)
    16 = @(1, 7)
    17 = @(1, 5)
    18 = @(1, 5)
    19 = @(1, 6)
    20 = @(1, 2)
    21 = @(1, 5)
    22 = @(1, 5)
    23 = @(1, 3)
    24 = @(1, 3)
    25 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
